$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists (Conta, Nome, Saldo) rows ordered by descending Saldo.
# Row 3 (account 004265173 / JULIA) needs to move down to row 7, with its
# Saldo updated to 9000; rows 4-7 (MERG, F, CATARINE, LUANA) each shift up
# one row to fill the gap. Use Range.Copy (not a literal .Value= write) so
# the numeric-looking "Conta" text (e.g. "004265173") keeps its original
# text type/leading zeros instead of being auto-coerced to a number.

# Stash row 3's account + name off to the side before it gets overwritten.
$ws.Range("A3:B3").Copy($ws.Range("Z1:AA1"))

# Shift rows 4-7 up into rows 3-6 (covers Conta, Nome, Saldo columns).
$ws.Range("A4:C7").Copy($ws.Range("A3:C6"))

# Drop the stashed JULIA account + name into row 7, then set her new Saldo.
$ws.Range("Z1:AA1").Copy($ws.Range("A7:B7"))
$ws.Cells.Item(7, 3).Value = 9000

# Clean up the scratch cells so they don't linger in the used range.
$ws.Range("Z1:AA1").Clear()
